# Updated cryptos list on Wed Sep  4 21:56:35 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to hold a literal text value (not auto-converted to a
# number/percentage by Excel's smart typing), while leaving the cell's style
# exactly as it was before (no lingering text/@ number format).
function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "57.903.05"
Set-TextValue $ws "E2" "  -0.14%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "2.420.09"
Set-TextValue $ws "E3" "  -1.27%  "

# Row 4 - TetherUSD
Set-TextValue $ws "E4" "  +0.15%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "511.60"
Set-TextValue $ws "E5" "  -2.48%  "

# Row 6 - Solana
Set-TextValue $ws "D6" "134.22"
Set-TextValue $ws "E6" "  +3.25%  "

# Row 7 - USDC
Set-TextValue $ws "D7" "0.998"
Set-TextValue $ws "E7" "  -0.17%  "

# Row 8 - XRP
Set-TextValue $ws "D8" "0.558"
Set-TextValue $ws "E8" "  -1.37%  "

# Row 9 - LidoStakedEther
Set-TextValue $ws "D9" "2.461.83"
Set-TextValue $ws "E9" "  +0.25%  "

# Row 10 - Dogecoin
Set-TextValue $ws "E10" "  +0.56%  "

# Row 11 - TRON
Set-TextValue $ws "E11" "  -0.54%  "

# Row 12 - Cardano
Set-TextValue $ws "D12" "0.325"
Set-TextValue $ws "E12" "  +0.55%  "

# Row 13 - Toncoin
Set-TextValue $ws "D13" "4.66"
Set-TextValue $ws "E13" "  -6.23%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws "D14" "2.896.25"
Set-TextValue $ws "E14" "  +0.34%  "

# Row 15 - WrappedBTC
Set-TextValue $ws "D15" "57.805.15"
Set-TextValue $ws "E15" "  -0.18%  "

# Row 16 - Avalanche
Set-TextValue $ws "D16" "22.03"
Set-TextValue $ws "E16" "  +2.18%  "

# Row 17 - ShibaInu
Set-TextValue $ws "E17" "  +1.45%  "

# Row 18 - WrappedEther
Set-TextValue $ws "D18" "2.484.56"
Set-TextValue $ws "E18" "  +1.30%  "

# Row 19 - Chainlink
Set-TextValue $ws "D19" "10.39"
Set-TextValue $ws "E19" "  -0.03%  "

# Row 20 - Polkadot
Set-TextValue $ws "D20" "4.16"
Set-TextValue $ws "E20" "  +0.97%  "

# Row 21 - BitcoinCash
Set-TextValue $ws "D21" "315.15"
Set-TextValue $ws "E21" "  +0.98%  "

# Row 22 - Uniswap
Set-TextValue $ws "D22" "6.47"
Set-TextValue $ws "E22" "  +5.31%  "

# Row 23 - Dai
Set-TextValue $ws "E23" "  +0.23%  "

# Row 24 - LEO
Set-TextValue $ws "E24" "  -1.23%  "

# Row 25 - Litecoin
Set-TextValue $ws "D25" "65.36"
Set-TextValue $ws "E25" "  +0.76%  "

# Row 26 - Binance-PegBSC-USD
Set-TextValue $ws "D26" "1.00"
Set-TextValue $ws "E26" "  -0.27%  "

# Row 27 - Kaspa
Set-TextValue $ws "E27" "  -0.30%  "

# Row 28 - Polygon
Set-TextValue $ws "D28" "0.383"
Set-TextValue $ws "E28" "  -4.66%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue $ws "D29" "7.62"
Set-TextValue $ws "E29" "  +5.16%  "

# Row 30 - Monero
Set-TextValue $ws "D30" "172.69"
Set-TextValue $ws "E30" "  -1.21%  "

# Row 31 - PEPE
Set-TextValue $ws "D31" "0.0₃0738"
Set-TextValue $ws "E31" "  +0.26%  "

# Row 32 - PancakeSwap
Set-TextValue $ws "D32" "1.70"
Set-TextValue $ws "E32" "  +0.38%  "

# Row 33 - Aptos
Set-TextValue $ws "E33" "  +0.01%  "

# Row 34 - Fetch.AI
Set-TextValue $ws "E34" "  +0.24%  "

# Row 35 - USDe (unchanged, per diff)

# Row 36 - FirstDigitalUSD
Set-TextValue $ws "D36" "0.993"
Set-TextValue $ws "E36" "  -0.42%  "

# Row 37 - EthereumClassic
Set-TextValue $ws "E37" "  +1.69%  "

# Row 38 - ImmutableX
Set-TextValue $ws "D38" "1.25"
Set-TextValue $ws "E38" "  +5.35%  "

# Row 39 - NEARProtocol
Set-TextValue $ws "D39" "3.89"
Set-TextValue $ws "E39" "  +2.72%  "

# Row 40 - OKB
Set-TextValue $ws "D40" "36.76"
Set-TextValue $ws "E40" "  +1.18%  "

# Row 41 - Stacks
Set-TextValue $ws "E41" "  +1.63%  "

# Row 42 - SuiNetwork
Set-TextValue $ws "D42" "0.813"
Set-TextValue $ws "E42" "  +0.46%  "

# Row 43 - Aave
Set-TextValue $ws "D43" "137.47"
Set-TextValue $ws "E43" "  +9.52%  "

# Row 44 - Filecoin
Set-TextValue $ws "D44" "3.43"
Set-TextValue $ws "E44" "  +1.30%  "

# Row 45 - RenderToken
Set-TextValue $ws "D45" "4.97"
Set-TextValue $ws "E45" "  +3.59%  "

# Rows 46 & 47 - Mantle/Bittensor swap places (with refreshed price/volume)
Set-TextValue $ws "B46" "Bittensor"
Set-TextValue $ws "C46" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws "D46" "257.84"
Set-TextValue $ws "E46" "  -0.53%  "

Set-TextValue $ws "B47" "Mantle"
Set-TextValue $ws "C47" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D47" "0.578"
Set-TextValue $ws "E47" "  -0.98%  "

# Row 48 - Stellar
Set-TextValue $ws "D48" "0.0922"
Set-TextValue $ws "E48" "  -0.06%  "

# Row 49 - Hedera
Set-TextValue $ws "D49" "0.0495"
Set-TextValue $ws "E49" "  +0.73%  "

# Row 50 - VeChain
Set-TextValue $ws "E50" "  +2.19%  "

# Row 51 - EnergySwap
Set-TextValue $ws "D51" "17.24"
Set-TextValue $ws "E51" "  +1.15%  "
